# Apply edits to the "Week 3" worksheet of the Time recording log workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 3")
$ws.Activate()

# Row 12: fill in the Stop time, mark the "x" (Class) column and the Lecture number.
$ws.Range("D12").Value = 0.97569444444444453
$ws.Range("I12").Value = "x"
$ws.Range("J12").Value = 35

# Row 13: fill in the date for the next log entry.
$ws.Range("B13").Value = 43514

# The long "Prep." comment block used to span rows 8-13; it now only spans 8-12,
# so row 13's Activity/Comments cells need to be un-merged and restored to the
# regular (fully-bordered) cell formatting used elsewhere in the table.
$ws.Range("G8:G13").UnMerge()
$ws.Range("H8:H13").UnMerge()
$ws.Range("G8:G12").Merge()
$ws.Range("H8:H12").Merge()
$ws.Range("G7:H7").Copy()
$ws.Range("G13:H13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the cursor/selection position as last saved in the file.
$ws.Range("D26").Select()
